$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("AddProduct")

# The "category" values in column B ("Cake") were renamed to "Sweet Cake"
# for all three product rows.
$ws.Range("B2").Value = "Sweet Cake"
$ws.Range("B3").Value = "Sweet Cake"
$ws.Range("B4").Value = "Sweet Cake"

# Match the selection left behind in the saved file (cell B4 selected).
$ws.Range("B4").Select()
